# Splits the single "Sales Data" sheet into two month-named sheets:
#   - "May 2024"   (existing rows, unchanged data, just re-homed/renamed)
#   - "April 2024" (new sheet, same customers/products one month earlier)
#
# The original sheet carried explicit custom column widths; the rebuilt
# "May 2024" sheet (and the brand-new "April 2024" sheet) should not, so we
# recreate the sheet from scratch rather than just renaming it in place.

$wb = $excel.ActiveWorkbook

$headers = @(
    "Customer Name",
    "Email",
    "Phone",
    "Product Purchased",
    "Quantity Purchased",
    "Total Sale Amount",
    "Date of Purchase"
)

# Flattened row data (7 values per row): Name, Email, Phone, Product, Qty, Total, DateSerial
$mayData = @(
    "John Smith","john.smith@example.com","555-1234","Laptop",2,2000,45413,
    "Emily Johnson","emily.johnson@example.com","555-5678","Smartphone",1,800,45415,
    "Michael Brown","michael.brown@example.com","555-9012","Tablet",3,1500,45417,
    "Sarah Davis","sarah.davis@example.com","555-3456","Smartwatch",1,300,45419,
    "David Wilson","david.wilson@example.com","555-7890","Headphones",2,150,45421,
    "Jennifer Martinez","jennifer.martinez@example.com","555-2345","Camera",1,600,45423,
    "Robert Garcia","robert.garcia@example.com","555-6789","Printer",1,400,45425,
    "Lisa Rodriguez","lisa.rodriguez@example.com","555-1234","External Hard Drive",2,200,45427,
    "Daniel Hernandez","daniel.hernandez@example.com","555-5678","Monitor",1,400,45429,
    "Maria Lopez","maria.lopez@example.com","555-9012","Keyboard",1,50,45431
)

$aprilData = @(
    "John Smith","john.smith@example.com","555-1234","Laptop",2,2000,45383,
    "Emily Johnson","emily.johnson@example.com","555-5678","Smartphone",1,800,45385,
    "Michael Brown","michael.brown@example.com","555-9012","Tablet",3,1500,45387,
    "Sarah Davis","sarah.davis@example.com","555-3456","Smartwatch",1,300,45389,
    "David Wilson","david.wilson@example.com","555-7890","Headphones",2,150,45391,
    "Jennifer Martinez","jennifer.martinez@example.com","555-2345","Camera",1,600,45393,
    "Robert Garcia","robert.garcia@example.com","555-6789","Printer",1,400,45395,
    "Lisa Rodriguez","lisa.rodriguez@example.com","555-1234","External Hard Drive",2,200,45397,
    "Daniel Hernandez","daniel.hernandez@example.com","555-5678","Monitor",1,400,45399,
    "Maria Lopez","maria.lopez@example.com","555-9012","Keyboard",1,50,45401
)

function Fill-Sheet($ws, $headers, $flatData) {
    for ($c = 0; $c -lt 7; $c++) {
        $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
    }
    $ws.Range("A1:G1").Font.Bold = $true

    $r = 2
    for ($i = 0; $i -lt $flatData.Count; $i += 7) {
        for ($c = 0; $c -lt 7; $c++) {
            $ws.Cells.Item($r, $c + 1).Value = $flatData[$i + $c]
        }
        $ws.Cells.Item($r, 7).NumberFormat = "yyyy-mm-dd h:mm:ss"
        $r++
    }
}

# 1) Build the replacement "May 2024" sheet ahead of the original sheet,
#    populate it, then drop the original so the rebuilt sheet carries no
#    leftover custom column widths.
$mayNew = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$mayNew.Name = "May2024Temp"
Fill-Sheet $wb.Worksheets.Item("May2024Temp") $headers $mayData

$wb.Worksheets.Item("Sales Data").Activate()
$wb.ActiveSheet.Delete() | Out-Null

$wb.Worksheets.Item("May2024Temp").Name = "May 2024"

# 2) Add the new "April 2024" sheet right after "May 2024" and populate it.
$aprilNew = $wb.Worksheets.Add($null, $wb.Worksheets.Item("May 2024"))
$aprilNew.Name = "April 2024"
Fill-Sheet $wb.Worksheets.Item("April 2024") $headers $aprilData

$wb.Worksheets.Item("May 2024").Activate()

Write-Host "Sheets now:" $wb.Worksheets.Count
